$wb = $excel.ActiveWorkbook

# --- Sheet "Annotations" (sheet1) becomes the active tab/sheet (was "Notes") ---
$wsAnn = $wb.Worksheets.Item("Annotations")
$wsNotes = $wb.Worksheets.Item("Notes")

# Insert two new rows before row 16 (format-from-below causes them to inherit
# the same look as the rows that get pushed down).
$wsAnn.Rows("16:17").Insert()

# Fill the new rows. Values are entered in this specific order so that the
# new shared-string entries end up appended in the expected order:
#   ET RoR value, BT RoR value, ~R2, ~R1
$wsAnn.Range("B16").Value = "ET RoR value"
$wsAnn.Range("B17").Value = "BT RoR value"
$wsAnn.Range("A17").Value = "~R2"
$wsAnn.Range("A16").Value = "~R1"

$wsAnn.Range("C16").Value = 10
$wsAnn.Range("C17").Value = 12

# Re-apply the "quote prefix" number style (same style used by row 19 / old
# row 17) to C16:C17 - entering plain numeric values resets formatting, so
# this must happen after the values are set.
$wsAnn.Range("C19").Copy()
$wsAnn.Range("C16:C17").PasteSpecial(-4122)

# Apply the same style to the existing numeric/text cells in column C rows 3-12
$wsAnn.Range("C15").Copy()
$wsAnn.Range("C3:C12").PasteSpecial(-4122)

# "Notes" sheet is no longer the active tab - select A11 there without activating it
# (do this before activating "Annotations" below, since selecting a range on a
# non-active sheet switches focus to that sheet)
$wsNotes.Range("A11").Select() | Out-Null

# Select D15 and make "Annotations" the active sheet/tab
$wsAnn.Activate() | Out-Null
$wsAnn.Range("D15").Select() | Out-Null
